$d = $word.ActiveDocument

$replacements = @(
    @{old = "## # A tibble: 66 x 3"; new = "## # A tibble: 66 x 2"},
    @{old = "##    ``Scientific name``            ``Common name``                 TotalPartners"; new = "##    ``Common name``                             TotalPartners"},
    @{old = "##    <chr>                        <chr>                                 <dbl>"; new = "##    <chr>                                             <dbl>"},
    @{old = "##  1 Calochortus persistens       SISKIYOU MARIPOSA LILY                    2"; new = "##  1 SISKIYOU MARIPOSA LILY                                2"},
    @{old = "##  2 Allium gooddingii            GOODING'S ONION                           1"; new = "##  2 GOODING'S ONION                                       1"},
    @{old = "##  3 Astragalus cusickii var. pa~ PACKARDS MILKVETCH                        1"; new = "##  3 PACKARDS MILKVETCH                                    1"},
    @{old = "##  4 Thymallus arcticus           ARCTIC GRAYLING- UPPER MISSO~             4"; new = "##  4 ARCTIC GRAYLING- UPPER MISSOURI RIVER DPS             4"},
    @{old = "##  5 Cimicifuga arizonica         ARIZONA BUGBANE                           2"; new = "##  5 ARIZONA BUGBANE                                       2"},
    @{old = "##  6 Lupinus aridus ssp. ashland~ ASHLAND LUPINE                            2"; new = "##  6 ASHLAND LUPINE                                        2"},
    @{old = "##  7 Pseudanophthalmus major      BEAVER CAVE BEETLE                        7"; new = "##  7 BEAVER CAVE BEETLE                                    7"},
    @{old = "##  8 Opuntia X multigeniculata    BLUE DIAMOND CHOLLA                       1"; new = "##  8 BLUE DIAMOND CHOLLA                                   1"},
    @{old = "##  9 Phacelia stellaris           BRAND'S PHACELIA                         20"; new = "##  9 BRAND'S PHACELIA                                     20"},
    @{old = "## 10 Fallicambarus gordoni        CAMP SHELBY BURROWING CRAYFI~             4"; new = "## 10 CAMP SHELBY BURROWING CRAYFISH                        4"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $found = $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $r.new
    } else {
        Write-Output "NOT FOUND: $($r.old)"
    }
}
